$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Existing rows 2:4 (Stimuli/210.jpg, 218.jpg, 222.jpg) are unchanged.
# Append the rest of the stimuli list (images that were deleted from the
# stimulus folder) below the existing three rows, sorted like the source.
$newValues = @(
    "Stimuli/228.jpg",
    "Stimuli/246.jpg",
    "Stimuli/251.jpg",
    "Stimuli/2703.jpg",
    "Stimuli/3051.jpg",
    "Stimuli/3160.jpg",
    "Stimuli/3185.jpg",
    "Stimuli/3301.jpg",
    "Stimuli/6562.jpg",
    "Stimuli/9031.jpg",
    "Stimuli/9040.jpg",
    "Stimuli/9042.jpg",
    "Stimuli/9043.jpg",
    "Stimuli/9145.jpg",
    "Stimuli/9160.jpg",
    "Stimuli/9184.jpg",
    "Stimuli/9904.jpg"
)

$row = 5
foreach ($val in $newValues) {
    $ws.Cells.Item($row, 1).Value = $val
    $row = $row + 1
}

# Update the selection to match the new data extent (A2:A21) and make the
# active cell A2, mirroring the saved view state after the edit.
$ws.Range("A2:A21").Select()
